# Bill Russell per-game averages: reorder the 13 season rows (2-14) into
# the new sequence captured in the target workbook, and switch the
# player_id (column C) from 415 to 413 for every data row.
#
# Approach: stage the current rows 2-14 into a scratch area far below the
# used range (rows 500-512), then paste them back into rows 2-14 in the
# new order by reading from the staged copies (whole-row copy/paste keeps
# cell types/formatting intact, unlike re-typing values through .Value).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Stage offset: row N (2..14) is staged at row (N + 498), i.e. 500..512.
$stageOffset = 498

for ($r = 2; $r -le 14; $r++) {
    $src = "A" + $r + ":AM" + $r
    $dstRow = $r + $stageOffset
    $dst = "A" + $dstRow + ":AM" + $dstRow
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial() | Out-Null
}

# Map: destination row (2..14) -> original row number (now staged).
$rowMap = @{
    2  = 14
    3  = 8
    4  = 6
    5  = 3
    6  = 2
    7  = 7
    8  = 9
    9  = 11
    10 = 10
    11 = 5
    12 = 4
    13 = 12
    14 = 13
}

foreach ($destRow in 2..14) {
    $origRow = $rowMap[$destRow]
    $stagedRow = $origRow + $stageOffset
    $src = "A" + $stagedRow + ":AM" + $stagedRow
    $dst = "A" + $destRow + ":AM" + $destRow
    $ws.Range($src).Copy() | Out-Null
    $ws.Range($dst).PasteSpecial() | Out-Null
}

# Clear the scratch area now that the reorder is complete.
$ws.Range("A500:AM512").Clear() | Out-Null

# player_id (column C) changes from 415 to 413 on every data row.
for ($r = 2; $r -le 14; $r++) {
    $ws.Cells.Item($r, 3).Value = 413
}

Write-Output "done"
